# Set cell A1 on the active sheet to the text "Some text".
# This causes Excel to create xl/sharedStrings.xml with the single
# shared string "Some text" and to reference it (type "s", shared
# string index 0) from the <c r="A1"> cell added to <sheetData>.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Some text"
